# Update column headers (C1, D1, E1) and remove F1 (scenario) by deleting column F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "immunity_all"
$ws.Range("D1").Value = "immunity_by_infection"
$ws.Range("E1").Value = "immunity_by_vaccination"

$ws.Range("D2").Value = 27.69389020031914
$ws.Range("E2").Value = 44.81348600949811
$ws.Range("D3").Value = 34.72731820352275
$ws.Range("E3").Value = 30.0654685149135
$ws.Range("D4").Value = 8.385509838998209
$ws.Range("E4").Value = 50.77993241202797
$ws.Range("D5").Value = 34.59715932609676
$ws.Range("E5").Value = 27.92224978552357
$ws.Range("D6").Value = 25.77179363093073
$ws.Range("E6").Value = 40.79653466175498
$ws.Range("D7").Value = 25.48576751707107
$ws.Range("E7").Value = 46.09740784647847
$ws.Range("D8").Value = 34.89980169084647
$ws.Range("E8").Value = 43.17335152827879
$ws.Range("D9").Value = 36.51663094387552
$ws.Range("E9").Value = 41.2609406211385
$ws.Range("D10").Value = 25.78915980298924
$ws.Range("E10").Value = 33.06171551933834
$ws.Range("D11").Value = 18.20733603914577
$ws.Range("E11").Value = 49.06748884478831
$ws.Range("D12").Value = 19.40399782244024
$ws.Range("E12").Value = 46.7876473795855
$ws.Range("D13").Value = 27.70215734825373
$ws.Range("E13").Value = 31.11422338968807
$ws.Range("D14").Value = 19.43413893175528
$ws.Range("E14").Value = 30.78172456414869
$ws.Range("D15").Value = 31.37646877571649
$ws.Range("E15").Value = 32.18610812380206
$ws.Range("D16").Value = 11.66495337663947
$ws.Range("E16").Value = 32.0718929433087
$ws.Range("D17").Value = 12.4734104754374
$ws.Range("E17").Value = 40.6322624770084
$ws.Range("D18").Value = 21.70286070832712
$ws.Range("E18").Value = 37.99377977496323
$ws.Range("D19").Value = 30.16151049625446
$ws.Range("E19").Value = 33.77182577149271
$ws.Range("D20").Value = 19.40922628354292
$ws.Range("E20").Value = 49.45006792033143
$ws.Range("D21").Value = 21.45232221387965
$ws.Range("E21").Value = 30.78858754939254
$ws.Range("D22").Value = 45.83626630870725
$ws.Range("E22").Value = 46.92851349613549
$ws.Range("D23").Value = 26.7045214971398
$ws.Range("E23").Value = 34.90912819560169
$ws.Range("D24").Value = 36.53138457110415
$ws.Range("E24").Value = 34.63768001177167
$ws.Range("D25").Value = 41.98198890358411
$ws.Range("E25").Value = 28.70112580140914
$ws.Range("D26").Value = 19.21369456079821
$ws.Range("E26").Value = 33.21704261938344
$ws.Range("D27").Value = 14.83640780524465
$ws.Range("E27").Value = 21.20765562758635
$ws.Range("D28").Value = 9.920663867483272
$ws.Range("E28").Value = 35.84053573000205
$ws.Range("D29").Value = 4.672043764308133
$ws.Range("E29").Value = 63.83578632724505
$ws.Range("D30").Value = 18.60505787800943
$ws.Range("E30").Value = 31.00558813095128
$ws.Range("D31").Value = 14.92697614050289
$ws.Range("E31").Value = 39.03882325982024
$ws.Range("D32").Value = 37.83080758338104
$ws.Range("E32").Value = 28.58253284614246
$ws.Range("D33").Value = 11.79075407545688
$ws.Range("E33").Value = 46.46831438225764
$ws.Range("D34").Value = 33.31268336586677
$ws.Range("E34").Value = 32.06924299347569
$ws.Range("D35").Value = 17.59085275656659
$ws.Range("E35").Value = 41.07807180373145
$ws.Range("D36").Value = 12.75057783300452
$ws.Range("E36").Value = 32.12846477458199
$ws.Range("D37").Value = 31.15217803890725
$ws.Range("E37").Value = 31.79724460889633
$ws.Range("D38").Value = 18.65124830140416
$ws.Range("E38").Value = 31.97828762071892
$ws.Range("D39").Value = 36.16080710921468
$ws.Range("E39").Value = 32.0937866595111
$ws.Range("D40").Value = 22.56036843415483
$ws.Range("E40").Value = 43.89165444035245
$ws.Range("D41").Value = 35.10966283751906
$ws.Range("E41").Value = 36.9604797000524
$ws.Range("D42").Value = 21.30486001453744
$ws.Range("E42").Value = 41.35631955539282
$ws.Range("D43").Value = 34.86257832809895
$ws.Range("E43").Value = 33.9806992108897
$ws.Range("D44").Value = 21.16674100224519
$ws.Range("E44").Value = 26.60163104363588
$ws.Range("D45").Value = 32.25444727090838
$ws.Range("E45").Value = 47.74323695337087
$ws.Range("D46").Value = 24.23237914435604
$ws.Range("E46").Value = 47.97383261996821
$ws.Range("D47").Value = 41.64584374479427
$ws.Range("E47").Value = 34.77326636654423
$ws.Range("D48").Value = 16.537385182788
$ws.Range("E48").Value = 21.76915728420372
$ws.Range("D49").Value = 32.46322064660833
$ws.Range("E49").Value = 53.76752517342644
$ws.Range("D50").Value = 20.59602260128642
$ws.Range("E50").Value = 30.27015658595338
$ws.Range("D51").Value = 20.78130542427902
$ws.Range("E51").Value = 39.91597965045125
$ws.Range("D52").Value = 18.83320113624375
$ws.Range("E52").Value = 31.39514712034604
$ws.Range("D53").Value = 51.3141226086817
$ws.Range("E53").Value = 39.69989399621305
$ws.Range("D54").Value = 21.1127970258937
$ws.Range("E54").Value = 36.82993955801849
$ws.Range("D55").Value = 35.5096309709038
$ws.Range("E55").Value = 36.50780809267105
$ws.Range("D56").Value = 19.11040153180073
$ws.Range("E56").Value = 41.34161779800483
$ws.Range("D57").Value = 26.95721536384507
$ws.Range("E57").Value = 52.37308612056123
$ws.Range("D58").Value = 31.48625875264324
$ws.Range("E58").Value = 46.89408129183581
$ws.Range("D59").Value = 30.97257244047038
$ws.Range("E59").Value = 36.16728376256484
$ws.Range("D60").Value = 24.49183690477609
$ws.Range("E60").Value = 37.86595008213377
$ws.Range("D61").Value = 14.03189750987513
$ws.Range("E61").Value = 36.27618982410679
$ws.Range("D62").Value = 18.4566166970859
$ws.Range("E62").Value = 37.81373061144841
$ws.Range("D63").Value = 54.75956598809967
$ws.Range("E63").Value = 29.95244885123448
$ws.Range("D64").Value = 30.83437827615269
$ws.Range("E64").Value = 43.88615948623971
$ws.Range("D65").Value = 32.72278078191085
$ws.Range("E65").Value = 35.69699860637994
$ws.Range("D66").Value = 11.87907372922596
$ws.Range("E66").Value = 43.12339610649641
$ws.Range("D67").Value = 47.82142554126911
$ws.Range("E67").Value = 37.12486576072786
$ws.Range("D68").Value = 11.98786859245466
$ws.Range("E68").Value = 18.12761394193073
$ws.Range("D69").Value = 10.8345354002511
$ws.Range("E69").Value = 64.79473456966376
$ws.Range("D70").Value = 16.95808325207584
$ws.Range("E70").Value = 45.92739951885694
$ws.Range("D71").Value = 34.26244252280595
$ws.Range("E71").Value = 37.01739049330526
$ws.Range("D72").Value = 17.83014746805538
$ws.Range("E72").Value = 39.48197512902068
$ws.Range("D73").Value = 13.74649178074346
$ws.Range("E73").Value = 36.57940148876396
$ws.Range("D74").Value = 28.15331827912842
$ws.Range("E74").Value = 40.13299380995644
$ws.Range("D75").Value = 8.622068442841503
$ws.Range("E75").Value = 39.23294401805383
$ws.Range("D76").Value = 23.20003292907899
$ws.Range("E76").Value = 37.39774150673931
$ws.Range("D77").Value = 25.15406867060747
$ws.Range("E77").Value = 30.53125203004993
$ws.Range("D78").Value = 34.19031354227033
$ws.Range("E78").Value = 31.46269485028784
$ws.Range("D79").Value = 31.79118989550678
$ws.Range("E79").Value = 20.84422268367878
$ws.Range("D80").Value = 31.92564871144535
$ws.Range("E80").Value = 37.82348438543815
$ws.Range("D81").Value = 34.10371916957992
$ws.Range("E81").Value = 23.94870609522823
$ws.Range("D82").Value = 51.97583239295578
$ws.Range("E82").Value = 27.98096525845754
$ws.Range("D83").Value = 27.91245198754862
$ws.Range("E83").Value = 33.39447632148848
$ws.Range("D84").Value = 38.59750695193734
$ws.Range("E84").Value = 30.80185770807167
$ws.Range("D85").Value = 35.47144548638347
$ws.Range("E85").Value = 27.54284935920796
$ws.Range("D86").Value = 26.88384331424021
$ws.Range("E86").Value = 29.00844595841327
$ws.Range("D87").Value = 36.33942303684284
$ws.Range("E87").Value = 36.28217008326812
$ws.Range("D88").Value = 25.89276081562196
$ws.Range("E88").Value = 32.14653777985414
$ws.Range("D89").Value = 17.76461880088823
$ws.Range("E89").Value = 49.22071810619293
$ws.Range("D90").Value = 26.29848783694937
$ws.Range("E90").Value = 38.59024854812698
$ws.Range("D91").Value = 15.11474723248645
$ws.Range("E91").Value = 38.83971799583165
$ws.Range("D92").Value = 31.47552777715326
$ws.Range("E92").Value = 42.64900603502772
$ws.Range("D93").Value = 10.77850834510538
$ws.Range("E93").Value = 48.99556026978128
$ws.Range("D94").Value = 20.87569682302821
$ws.Range("E94").Value = 43.86128860418635
$ws.Range("D95").Value = 42.58449150245828
$ws.Range("E95").Value = 37.90085510855813
$ws.Range("D96").Value = 9.523476068765097
$ws.Range("E96").Value = 53.58060789650659
$ws.Range("D97").Value = 31.98827013782962
$ws.Range("E97").Value = 26.67571174349406
$ws.Range("D98").Value = 26.51396420839526
$ws.Range("E98").Value = 27.64552447644686
$ws.Range("D99").Value = 34.3832504451403
$ws.Range("E99").Value = 36.65932971747905
$ws.Range("D100").Value = 22.50289469054428
$ws.Range("E100").Value = 30.7038529202802
$ws.Range("D101").Value = 25.98752598752599
$ws.Range("E101").Value = 39.77196842482534

# Remove the old "scenario" column (F) entirely; Excel shifts remaining
# columns left and the used range / dimension shrinks to A1:E101.
$ws.Range("F1:F101").Delete()
